# Commit: "Comment => Annotation in XLSX"
# Rename the header cell text in Sheet1!H1 from "Comment" to "Annotation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Annotation"
